# The sheet tracks weekly Berenjena (eggplant) price observations for the
# "Terminal La Palmera de La Serena" market. This edit inserts one new
# observation row at sheet row 165 (pushing the existing rows 165-253 down
# to 166-254), extending the used range from A1:R253 to A1:R254.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 165 - shifts rows 165..253 down to 166..254
# and carries the date-column number format (style) down with them.
$ws.Rows(165).Insert()

# Populate the newly inserted row 165 with the new observation. The
# "constant" columns (A,B,C,E,F,G,H,I,O,R) carry the same values used by
# every other row on this sheet.
$ws.Range("A165").Value2 = 8
$ws.Range("B165").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C165").Value2 = "Coquimbo"
$ws.Range("D165").Value2 = 45097
$ws.Range("E165").Value2 = 4
$ws.Range("F165").Value2 = 100112001
$ws.Range("G165").Value2 = "Berenjena"
$ws.Range("H165").Value2 = "Sin especificar"
$ws.Range("I165").Value2 = "Primera"
$ws.Range("J165").Value2 = 400
$ws.Range("K165").Value2 = 8000
$ws.Range("L165").Value2 = 8500
$ws.Range("M165").Value2 = 8250
$ws.Range("N165").Value2 = "`$/caja 50 unidades"
$ws.Range("O165").Value2 = "Regi$([char]0xF3)n de Arica y Parinacota"
$ws.Range("P165").Value2 = 165
$ws.Range("Q165").Value2 = 50
$ws.Range("R165").Value2 = "Hortaliza"

Write-Host "Inserted new row 165 for Berenjena weekly log; used range now A1:R254"
